$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ciudad = "Bucaramanga - Floridablanca - Piedecuesta"

$rows = @(
    @{ Row = 2;  Manzana = "'6800110000000002800616"; Lat = 7.14499230679;  Lon = -73.124413 },
    @{ Row = 3;  Manzana = "'6800110000000002850125"; Lat = 7.15298783395;  Lon = -73.12597473069999 },
    @{ Row = 4;  Manzana = "'6800110000000002850513"; Lat = 7.14997382542;  Lon = -73.1239001484 },
    @{ Row = 5;  Manzana = "'6800110000000009440321"; Lat = 7.09960313404;  Lon = -73.10988221700001 },
    @{ Row = 6;  Manzana = "'6827610000000000090315"; Lat = 7.06004692805;  Lon = -73.0866898312 },
    @{ Row = 7;  Manzana = "'6827610000000000030630"; Lat = 7.07496927963;  Lon = -73.0833930158 },
    @{ Row = 8;  Manzana = "'6800110000000012380513"; Lat = 7.11615102554;  Lon = -73.1060115377 },
    @{ Row = 9;  Manzana = "'6827610000000000110231"; Lat = 7.08236877722;  Lon = -73.109770799 },
    @{ Row = 10; Manzana = "'6800110000000010620105"; Lat = 7.09316301106;  Lon = -73.1110351889 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $ciudad
    $ws.Cells.Item($row, 3).Value = $r.Manzana
    $ws.Cells.Item($row, 5).Value = $r.Lat
    $ws.Cells.Item($row, 6).Value = $r.Lon
}
